$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'311.93"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'0.86%"
$ws.Range('E2').ClearFormats()

$ws.Range('D3').Value = "'37.80"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "'0.23%"
$ws.Range('E3').ClearFormats()

$ws.Range('D4').Value = "'5.122"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "'0.16%"
$ws.Range('E4').ClearFormats()

$ws.Range('D5').Value = "'0.07910"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'0.67%"
$ws.Range('E5').ClearFormats()

$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Value = "'1.911"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'-3.01%"
$ws.Range('E6').ClearFormats()

$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = "'8.273"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'-0.29%"
$ws.Range('E7').ClearFormats()

$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = "'2.860"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'-9.33%"
$ws.Range('E8').ClearFormats()

$ws.Range('D9').Value = "'0.9305"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'0.27%"
$ws.Range('E9').ClearFormats()

$ws.Range('D10').Value = "'0.1227"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'-8.77%"
$ws.Range('E10').ClearFormats()

$ws.Range('D11').Value = "'0.1923"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'-3.71%"
$ws.Range('E11').ClearFormats()

$ws.Range('D12').Value = "'0.09100"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "'1.05%"
$ws.Range('E12').ClearFormats()

$ws.Range('D13').Value = "'0.03333"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "'-2.93%"
$ws.Range('E13').ClearFormats()

$ws.Range('D14').Value = "'0.09650"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'-0.83%"
$ws.Range('E14').ClearFormats()

$ws.Range('D15').Value = "'0.001387"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'-0.20%"
$ws.Range('E15').ClearFormats()

$ws.Range('D16').Value = "'0.005810"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'-2.30%"
$ws.Range('E16').ClearFormats()

$ws.Range('D17').Value = "'3.535"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'-1.37%"
$ws.Range('E17').ClearFormats()

$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = "'4.405"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'0.90%"
$ws.Range('E18').ClearFormats()

$ws.Range('D19').Value = "'0.3408"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'-1.69%"
$ws.Range('E19').ClearFormats()

$ws.Range('D20').Value = "'5.278"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "'5.35%"
$ws.Range('E20').ClearFormats()

$ws.Range('D21').Value = "'0.1278"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'-1.29%"
$ws.Range('E21').ClearFormats()

$ws.Range('D22').Value = "'0.2618"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'5.13%"
$ws.Range('E22').ClearFormats()

$ws.Range('D24').Value = "'0.04362"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'0.90%"
$ws.Range('E24').ClearFormats()

$ws.Range('D25').Value = "'0.001240"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "'1.70%"
$ws.Range('E25').ClearFormats()

$ws.Range('D26').Value = "'0.004306"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "'-5.36%"
$ws.Range('E26').ClearFormats()

$ws.Range('D27').Value = "'0.0001221"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'-9.69%"
$ws.Range('E27').ClearFormats()

$ws.Range('E39').Value = "'-6.37%"
$ws.Range('E39').ClearFormats()

$ws.Range('D40').Value = "'0.05180"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'2.67%"
$ws.Range('E40').ClearFormats()

$ws.Range('D41').Value = "'0.007661"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'3.66%"
$ws.Range('E41').ClearFormats()

$ws.Range('D42').Value = "'0.009145"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "'-7.85%"
$ws.Range('E42').ClearFormats()

$ws.Range('D43').Value = "'0.1363"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "'0.72%"
$ws.Range('E43').ClearFormats()

$ws.Range('D44').Value = "'0.002052"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'3.49%"
$ws.Range('E44').ClearFormats()

$ws.Range('D45').Value = "'0.008615"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'-1.74%"
$ws.Range('E45').ClearFormats()

$ws.Range('D46').Value = "'0.00006692"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'0.24%"
$ws.Range('E46').ClearFormats()

$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'-0.05%"
$ws.Range('E47').ClearFormats()

$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').Value = "'0.001201"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'-7.75%"
$ws.Range('E48').ClearFormats()

$ws.Range('B49').Value = 'BOLO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D49').Value = "'0.002987"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'-0.57%"
$ws.Range('E49').ClearFormats()

$ws.Range('D50').Value = "'0.00002102"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'-0.05%"
$ws.Range('E50').ClearFormats()

$ws.Range('D51').Value = "'0.0002002"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'-0.05%"
$ws.Range('E51').ClearFormats()
